# ---------------------------------------------------------------------------
# Applies the changes described by the diff:
#  1. Insert a new worksheet "vecka_8" right after "vecka_50" and populate it
#     with a weekly schedule (anstalldID / StartTid / SlutTid / datum / CreatedBy).
#  2. Update a couple of sheet-view selections ("vecka_50" and "Ordrar").
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "vecka_8" sheet right after "vecka_50"
# ---------------------------------------------------------------------------
$weekSrc   = $wb.Worksheets.Item("vecka_50")
$newSheet  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $weekSrc)
$newSheet.Name = "vecka_8"

# Header row
$newSheet.Range("A1").Value = "anstalldID"
$newSheet.Range("D1").Value = "datum"
$newSheet.Range("E1").Value = "CreatedBy"

# Data rows (anstalldID / datum / CreatedBy) -- written column by column like
# the source workbook (A, D, E all reuse existing shared strings)
$newSheet.Range("A2").Value = "75E36D07-0BC5-4687-A94D-3BD9174EF194"
$newSheet.Range("A3").Value = "75E36D07-0BC5-4687-A94D-3BD9174EF194"
$newSheet.Range("A4").Value = "75E36D07-0BC5-4687-A94D-3BD9174EF194"
$newSheet.Range("A5").Value = "75E36D07-0BC5-4687-A94D-3BD9174EF194"
$newSheet.Range("A6").Value = "75E36D07-0BC5-4687-A94D-3BD9174EF194"
$newSheet.Range("A7").Value = "75E36D07-0BC5-4687-A94D-3BD9174EF194"

$newSheet.Range("D2").Value = "2016-12-01 14:59:34.391"
$newSheet.Range("D3").Value = "2016-12-01 14:59:34.392"
$newSheet.Range("D4").Value = "2016-12-01 14:59:34.393"
$newSheet.Range("D5").Value = "2016-12-01 14:59:34.394"
$newSheet.Range("D6").Value = "2016-12-01 14:59:34.395"
$newSheet.Range("D7").Value = "2016-12-01 14:59:34.396"

$newSheet.Range("E2").Value = "A3B762F3-79F8-49B2-8722-354505C82FF4"
$newSheet.Range("E3").Value = "A3B762F3-79F8-49B2-8722-354505C82FF4"
$newSheet.Range("E4").Value = "A3B762F3-79F8-49B2-8722-354505C82FF4"
$newSheet.Range("E5").Value = "A3B762F3-79F8-49B2-8722-354505C82FF4"
$newSheet.Range("E6").Value = "A3B762F3-79F8-49B2-8722-354505C82FF4"
$newSheet.Range("E7").Value = "A3B762F3-79F8-49B2-8722-354505C82FF4"

# StartTid / SlutTid values -- this exact entry order reproduces the shared
# string table ordering seen in the target workbook.
$newSheet.Range("B2").Value = "2017-02-21  10:30:00.00"
$newSheet.Range("B3").Value = "2017-02-21  13:00:00.00"
$newSheet.Range("C2").Value = "2017-02-21  14:30:00.00"
$newSheet.Range("C3").Value = "2017-02-21  20:00:00.00"
$newSheet.Range("B4").Value = "2017-02-22  10:00:00.00"
$newSheet.Range("B5").Value = "2017-02-22  14:30:00.00"
$newSheet.Range("C4").Value = "2017-02-22  14:00:00.00"
$newSheet.Range("C5").Value = "2017-02-22  21:00:00.00"
$newSheet.Range("B6").Value = "2017-02-23  12:00:00.00"
$newSheet.Range("B7").Value = "2017-02-23  16:30:00.00"
$newSheet.Range("C6").Value = "2017-02-23  16:00:00.00"
$newSheet.Range("C7").Value = "2017-02-23  20:30:00.00"
$newSheet.Range("B1").Value = "StartTid"
$newSheet.Range("C1").Value = "SlutTid"

# Apply the same number format (built-in date/time, style index 2 in the
# original workbook) that the other "vecka_*" sheets use for their
# start/end-time columns, by copying the formatting from "vecka_50".
$weekSrc.Range("B2").Copy()
$newSheet.Range("B2:C7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Column widths roughly matching the other week sheets (best effort; exact
# "best fit" pixel widths cannot be reproduced outside real Excel rendering).
$newSheet.Columns.Item(1).ColumnWidth = 37.5924479166667
$newSheet.Columns.Item(2).ColumnWidth = 20.5924479166667
$newSheet.Columns.Item(3).ColumnWidth = 20.5924479166667
$newSheet.Columns.Item(4).ColumnWidth = 20.5924479166667
$newSheet.Columns.Item(5).ColumnWidth = 21.1666666666667
$newSheet.Columns.Item(6).ColumnWidth = 36.5924479166667

# Selection / view state for the new sheet
$newSheet.Range("A7").Select()

# ---------------------------------------------------------------------------
# 2. Update existing sheet selections
# ---------------------------------------------------------------------------
$weekSrc.Range("A20:E25").Select()

$ordrar = $wb.Worksheets.Item("Ordrar")
$ordrar.Range("L24").Select()

$ettFucking = $wb.Worksheets.Item("Ett fucking blad")
$ettFucking.Activate()
$ettFucking.Range("B4").Select()
